$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

# Wipe the old content & formatting of the Input sheet data range entirely
$ws.Range("A1:Q5").Clear()

# --- Header row (row 1), columns A:P, no special style (plain text headers) ---
$headers = @("발주일자","납기일자","거래처명","거래처 이메일","납품처명","납품처 이메일","프로젝트명","대분류","중분류","소분류","품목명","규격","수량","단가","총금액","비고")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value2 = $headers[$i]
}

# --- Data rows 2-5 ---
# Columns A and B hold date-like text; force text format so Excel doesn't convert them to date serials
$ws.Range("A2:B5").NumberFormat = "@"

$data = @(
    @("2025-08-31","2025-09-13","센트럴머시너리","센트럴머시너리@example.com","힐스테이트 도곡동1차","delivery@example.com","힐스테이트 도곡동1차","2. 부자재비","1) 판넬","기타","3월 절삭","KS규격-1",1,0,0),
    @("2025-09-07","2025-09-15","센트럴머시너리","센트럴머시너리@example.com","힐스테이트 도곡동1차","delivery@example.com","힐스테이트 도곡동1차","2. 부자재비","1) 판넬","기타","사각와샤 40*40*4T","KS규격-2",280,150,46200),
    @("2025-08-25","2025-09-25","센트럴머시너리","센트럴머시너리@example.com","힐스테이트 도곡동1차","delivery@example.com","힐스테이트 도곡동1차","2. 부자재비","1) 판넬","기타","I-01 120*50*5880*7T 9T","KS규격-3",44,221250,10708500),
    @("2025-09-09","2025-10-01","센트럴머시너리","센트럴머시너리@example.com","힐스테이트 도곡동1차","delivery@example.com","힐스테이트 도곡동1차","5. 운반비","일반자재","기타","3월 운반비","KS규격-4",1,0,0)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    $excelRow = $r + 2
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($excelRow, $c + 1).Value2 = $row[$c]
    }
}

# --- Sheets 갑지 and 을지: drop the leftover truly-empty "비고" cells (I2:I5) ---
$ws2 = $wb.Worksheets.Item("갑지")
$ws2.Range("I2:I5").ClearContents()

$ws3 = $wb.Worksheets.Item("을지")
$ws3.Range("I2:I5").ClearContents()
